# Automatische test-sync: 2025-08-03 18:19:50
# Adds a new log entry (row 32) to the "Logs" sheet and updates the
# "Inkoop / Bestellingen" tally on the "Dashboard" sheet from 4 to 5.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A32").Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("B32").Value = "mailmind.test@zohomail.eu"
$logs.Range("C32").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("D32").Value = "Inkoop / Bestellingen"
$logs.Range("E32").Value = "Beste klant,`nHartelijk dank voor je interesse in onze producten. Helaas kan ik je niet helpen met deze specifieke bestelling via e-mail. Voor het plaatsen van een bestelling verwijs ik je graag door naar onze website of klantenservice. Mocht je hulp nodig hebben bij het plaatsen van een bestelling, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F32").Value = "2025-08-03 18:19:13"
$logs.Range("G32").Value = "Ja"
$logs.Range("H32").Value = "Nee"
$logs.Range("I32").Value = "Ja"
$logs.Range("J32").Value = "Nee"

# Setting E32 with embedded newlines auto-expands the row height; restore
# the default (no explicit/custom height), matching the other rows.
$logs.Rows.Item(32).AutoFit()

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B5").Value = 5

# Conditional formatting ranges in "Logs" must grow to include the new row.
$ccols = @("D", "G", "H", "I", "J")
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + "31")
    $newRange = $logs.Range($col + "2:" + $col + "32")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
